$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Nganh" data rows, keep A1:B1 formatting to reuse on the new header
$ws.Range("A2:B4").ClearContents()

# New header row describing mon hoc (subject) columns
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "Tên môn học"
$ws.Range("C1").Value = "Tổng STC"
$ws.Range("D1").Value = "Lý thuyết"
$ws.Range("E1").Value = "Thực hành"
$ws.Range("F1").Value = "Bắt buộc"

# Re-use the existing header formatting (same highlighted/centered style as A1) for the new columns
$ws.Range("A1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Totals row with self-referencing SUM formulas, as in the source sheet
$ws.Range("C2").Formula = "=SUM(C2:C1)"
$ws.Range("D2").Formula = "=SUM(D2:D1)"
$ws.Range("E2").Formula = "=SUM(E2:E1)"

# Column B narrows to fit the new shorter header text
$ws.Columns("B").AutoFit()
